$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-08-26 Monday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-08-27 Tuesday", 2)

# Update the division-problem table cells by position, since several
# cells share identical original text but need different replacements.
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="59÷6=9, 5"},
    @{Row=1;  Col=2; Text="46÷9=5, 1"},
    @{Row=1;  Col=3; Text="75÷8=9, 3"},
    @{Row=1;  Col=4; Text="43÷4=10, 3"},
    @{Row=1;  Col=5; Text="93÷2=46, 1"},

    @{Row=5;  Col=1; Text="71÷5=14, 1"},
    @{Row=5;  Col=2; Text="58÷5=11, 3"},
    @{Row=5;  Col=3; Text="96÷6=16, 0"},
    @{Row=5;  Col=4; Text="13÷3=4, 1"},
    @{Row=5;  Col=5; Text="95÷4=23, 3"},

    @{Row=9;  Col=1; Text="33÷8=4, 1"},
    @{Row=9;  Col=2; Text="83÷5=16, 3"},
    @{Row=9;  Col=3; Text="16÷8=2, 0"},
    @{Row=9;  Col=4; Text="36÷7=5, 1"},
    @{Row=9;  Col=5; Text="59÷5=11, 4"},

    @{Row=13; Col=1; Text="43÷7=6, 1"},
    @{Row=13; Col=2; Text="25÷9=2, 7"},
    @{Row=13; Col=3; Text="87÷7=12, 3"},
    @{Row=13; Col=4; Text="72÷5=14, 2"},
    @{Row=13; Col=5; Text="92÷8=11, 4"},

    @{Row=17; Col=1; Text="32÷5=6, 2"},
    @{Row=17; Col=2; Text="48÷6=8, 0"},
    @{Row=17; Col=3; Text="23÷7=3, 2"},
    @{Row=17; Col=4; Text="30÷4=7, 2"},
    @{Row=17; Col=5; Text="84÷2=42, 0"}
)

foreach ($u in $updates) {
    $t.Cell($u.Row, $u.Col).Range.Text = $u.Text
}
